$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Data updates - shared test phone numbers / virtual account renumbered.
#    "Paket Data Purchase" sheet: the "phone number for paket data" test
#    value +6281252930385 -> +6281252930365, and its related virtual
#    account number 8000081252930385 -> 8000081252930365.
# ---------------------------------------------------------------------------
$wsPurchase = $wb.Worksheets.Item("Paket Data Purchase")
$purchaseRange = $wsPurchase.UsedRange
$purchaseRange.Replace("+6281252930385", "+6281252930365", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole) | Out-Null
$purchaseRange.Replace("8000081252930385", "8000081252930365", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole) | Out-Null

# "Paket Data History" sheet: the "valid phone number" positive-case value
# +6281252930393 -> +6281252930363 (used for both phoneNumber and the
# echoed responseBodyRequest).
$wsHistory = $wb.Worksheets.Item("Paket Data History")
$historyRange = $wsHistory.UsedRange
$historyRange.Replace("+6281252930393", "+6281252930363", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole) | Out-Null

# ---------------------------------------------------------------------------
# 2) Column width tweak on "Paket Data Purchase" - column B widened.
# ---------------------------------------------------------------------------
$wsPurchase.Columns.Item(2).ColumnWidth = 22.6667

# ---------------------------------------------------------------------------
# 3) View state - selected cell / scroll position per sheet.
# ---------------------------------------------------------------------------
$wsList = $wb.Worksheets.Item("Paket Data List")
$wsList.Activate()
$wsList.Range("A31").Select()

$wsPurchase.Activate()
$wsPurchase.Range("J76").Select()

$wsHistory.Activate()
$wsHistory.Range("E13").Select()

Write-Host "done"
